$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Update the existing note text in A14 to append additional context about energy use data.
$ws.Range("A14").Value = "sectors.  It also includes data on energy use per ton CO2 sequestered."

# Insert 7 new rows right after row 14 (before the existing blank separator row 15).
# This pushes the original blank row (and everything below it) down by 7, and the new
# rows inherit the formatting of row 14 (the paragraph style used in A12:A14).
$ws.Rows("15:21").Insert()

# The first of the inserted rows should actually be the (unchanged) blank separator row,
# so restore its formatting back to match the blank-row style used elsewhere (e.g. A11/A22).
$ws.Range("A15").Font.Bold = $true

# Fill in the new note text for rows 16-21 (these already carry the paragraph style from the insert).
$ws.Range("A16").Value = "The energy use value here should exclude any energy that is not additional"
$ws.Range("A17").Value = "to the energy already accounted for in indst/BIFUbC or the heat rates in the"
$ws.Range("A18").Value = "electricity sector.  For example, if a steel mill has excess process heat that"
$ws.Range("A19").Value = "isn't being used, and they use it to power CCS, that heat should be excluded"
$ws.Range("A20").Value = "from here, since it does not increase the overall energy demand of the steel"
$ws.Range("A21").Value = "mill."
